$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2094.5312
$ws.Range("I98").Value = 2090.2856
$ws.Range("K98").Value = 2090.2856
$ws.Range("M98").Value = -592.2856000000002

$ws.Range("H107").Value = 581.0526
$ws.Range("I107").Value = 340
$ws.Range("J107").Value = 1866.6666
$ws.Range("K107").Value = 340
$ws.Range("L107").Value = 1866.6666
$ws.Range("M107").Value = 1580
$ws.Range("N107").Value = -5706.6666

$ws.Range("H122").Value = 2094.5312
$ws.Range("I122").Value = 2090.2856
$ws.Range("K122").Value = 6270.8568
$ws.Range("M122").Value = -3820.8568

$ws.Range("H127").Value = 3419.6924
$ws.Range("J127").Value = 5400
$ws.Range("L127").Value = 16200
$ws.Range("N127").Value = -26120

$ws.Range("H131").Value = 2386.2727
$ws.Range("J131").Value = 4345.3
$ws.Range("L131").Value = 13035.9
$ws.Range("N131").Value = -23115.9

$ws.Range("H137").Value = 1167.8049
$ws.Range("I137").Value = 942.3913
$ws.Range("J137").Value = 1455.8334
$ws.Range("K137").Value = 2827.1739
$ws.Range("L137").Value = 4367.5002
$ws.Range("M137").Value = -277.1738999999998
$ws.Range("N137").Value = -9467.5002

$ws.Range("H138").Value = 1817.0448
$ws.Range("I138").Value = 1552.0962
$ws.Range("J138").Value = 2735.5334
$ws.Range("K138").Value = 4656.2886
$ws.Range("L138").Value = 8206.600199999999
$ws.Range("M138").Value = 483.7114000000001
$ws.Range("N138").Value = -18486.6002

$ws.Range("H139").Value = 51193.855
$ws.Range("J139").Value = 51193.855
$ws.Range("L139").Value = 51193.855
$ws.Range("N139").Value = -61473.855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3141.482
$ws.Range("I32").Value = 2661.4304
$ws.Range("J32").Value = 12622.5
$ws.Range("K32").Value = 2661.4304
$ws.Range("L32").Value = 12622.5
$ws.Range("M32").Value = -2374.4304
$ws.Range("N32").Value = -13196.5

$ws.Range("H61").Value = 3012
$ws.Range("I61").Value = 2171.158
$ws.Range("J61").Value = 11000
$ws.Range("K61").Value = 2171.158
$ws.Range("L61").Value = 11000
$ws.Range("M61").Value = -1959.158
$ws.Range("N61").Value = -11424

$ws.Range("H74").Value = 2031.5
$ws.Range("I74").Value = 1695.3077
$ws.Range("J74").Value = 2655.8572
$ws.Range("K74").Value = 1695.3077
$ws.Range("L74").Value = 2655.8572
$ws.Range("M74").Value = -821.3077000000001
$ws.Range("N74").Value = -4403.8572

$ws.Range("H77").Value = 2031.5
$ws.Range("I77").Value = 1695.3077
$ws.Range("J77").Value = 2655.8572
$ws.Range("K77").Value = 8476.538500000001
$ws.Range("L77").Value = 13279.286
$ws.Range("M77").Value = -4108.538500000001
$ws.Range("N77").Value = -22015.286

$ws.Range("H136").Value = 3012
$ws.Range("I136").Value = 2171.158
$ws.Range("J136").Value = 11000
$ws.Range("K136").Value = 6513.474
$ws.Range("L136").Value = 33000
$ws.Range("M136").Value = -3963.474
$ws.Range("N136").Value = -38100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1554
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H105").Value = 2128.3242
$ws.Range("I105").Value = 2095.4243
$ws.Range("K105").Value = 2095.4243
$ws.Range("M105").Value = -348.4243000000001

$ws.Range("H134").Value = 8171.1665
$ws.Range("I134").Value = 8380.0625
$ws.Range("J134").Value = 6500
$ws.Range("K134").Value = 25140.1875
$ws.Range("L134").Value = 19500
$ws.Range("M134").Value = -22605.1875
$ws.Range("N134").Value = -24570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 988
$ws.Range("J16").Value = 1999.5
$ws.Range("L16").Value = 1999.5
$ws.Range("N16").Value = -2573.5

$ws.Range("H22").Value = 870.4706
$ws.Range("J22").Value = 1512.25
$ws.Range("L22").Value = 1512.25
$ws.Range("N22").Value = -2212.25

$ws.Range("H31").Value = 1936.7241
$ws.Range("I31").Value = 1527.0588
$ws.Range("J31").Value = 2517.0833
$ws.Range("K31").Value = 1527.0588
$ws.Range("L31").Value = 2517.0833
$ws.Range("M31").Value = -1232.0588
$ws.Range("N31").Value = -3107.0833

$ws.Range("H34").Value = 1936.7241
$ws.Range("I34").Value = 1527.0588
$ws.Range("J34").Value = 2517.0833
$ws.Range("K34").Value = 1527.0588
$ws.Range("L34").Value = 2517.0833
$ws.Range("M34").Value = -1325.0588
$ws.Range("N34").Value = -2921.0833

$ws.Range("H58").Value = 2559387.5
$ws.Range("I58").Value = 3953820.5
$ws.Range("J58").Value = 2926.8333
$ws.Range("K58").Value = 3953820.5
$ws.Range("L58").Value = 2926.8333
$ws.Range("M58").Value = -3953617.5
$ws.Range("N58").Value = -3332.8333

$ws.Range("H99").Value = 1733.75
$ws.Range("I99").Value = 1733.75
$ws.Range("K99").Value = 1733.75
$ws.Range("M99").Value = -235.75

$ws.Range("H113").Value = 988
$ws.Range("J113").Value = 1999.5
$ws.Range("L113").Value = 1999.5
$ws.Range("N113").Value = -6339.5

$ws.Range("H126").Value = 1733.75
$ws.Range("I126").Value = 1733.75
$ws.Range("K126").Value = 5201.25
$ws.Range("M126").Value = -2731.25

$ws.Range("H134").Value = 1744.0638
$ws.Range("J134").Value = 2210.9
$ws.Range("L134").Value = 6632.700000000001
$ws.Range("N134").Value = -11702.7

$ws.Range("H136").Value = 2559387.5
$ws.Range("I136").Value = 3953820.5
$ws.Range("J136").Value = 2926.8333
$ws.Range("K136").Value = 11861461.5
$ws.Range("L136").Value = 8780.499899999999
$ws.Range("M136").Value = -11858911.5
$ws.Range("N136").Value = -13880.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 377
$ws.Range("I18").Value = 323.16666
$ws.Range("K18").Value = 969.4999799999999
$ws.Range("M18").Value = -800.4999799999999

$ws.Range("H97").Value = 126.666664
$ws.Range("I97").Value = 140
$ws.Range("K97").Value = 420
$ws.Range("M97").Value = 76

$ws.Range("H126").Value = 5999.75
$ws.Range("I126").Value = 5999
$ws.Range("K126").Value = 17997
$ws.Range("M126").Value = -13057

$ws.Range("H131").Value = 9179.766
$ws.Range("J131").Value = 9634.709999999999
$ws.Range("L131").Value = 28904.13
$ws.Range("N131").Value = -38984.13

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 727260.5
$ws.Range("I132").Value = 895677
$ws.Range("J132").Value = 3069.7
$ws.Range("K132").Value = 2687031
$ws.Range("L132").Value = 9209.099999999999
$ws.Range("M132").Value = -2684501
$ws.Range("N132").Value = -14269.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 11379.223
$ws.Range("I16").Value = 14477.857
$ws.Range("J16").Value = 534
$ws.Range("K16").Value = 14477.857
$ws.Range("L16").Value = 534
$ws.Range("M16").Value = -14307.857
$ws.Range("N16").Value = -874

$ws.Range("H61").Value = 2399
$ws.Range("I61").Value = 2276.5386
$ws.Range("K61").Value = 2276.5386
$ws.Range("M61").Value = -2074.5386

$ws.Range("H113").Value = 2399
$ws.Range("I113").Value = 2276.5386
$ws.Range("K113").Value = 2276.5386
$ws.Range("M113").Value = -106.5385999999999

$ws.Range("H136").Value = 2616.4546
$ws.Range("J136").Value = 5131.75
$ws.Range("L136").Value = 15395.25
$ws.Range("N136").Value = -20495.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1146.3448
$ws.Range("I132").Value = 830.8837
$ws.Range("K132").Value = 2492.6511
$ws.Range("M132").Value = 37.34889999999996
